# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the file
# 18a43bd9-021d-4ca6-82f4-2da40df85fea on both the zh-cn and de-de
# localization-status sheets, reflecting a newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 on each sheet corresponds to the 18a43bd9... file (column D = "Latest Handoff Datetime")
$wsZhCn.Range("D4").Value = "2016-03-08 06:14:37"
$wsDeDe.Range("D4").Value = "2016-03-08 06:14:46"
